# Actualización automática 2025-10-29 11:30:07
#
# Updates sales figures for "CASTRO ALCIVAR EDA MARIA" across the three
# worksheets of the workbook: the per-group monthly sales sheet, the
# monthly-by-column sales sheet, and the compliance (budget vs sales)
# summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# LAVABOS sale recorded for client in row 4 (was 0)
$wsGrupo.Range("I4").Value = 134.51

# PIEDRA SINTERIZADA / PORCELANATO sales increase for client in row 15
$wsGrupo.Range("L15").Value = 2732.74
$wsGrupo.Range("M15").Value = 2375.06

# PORCELANATO sale recorded for client in row 29 (was 0)
$wsGrupo.Range("M29").Value = 366.83

# LAVABOS sale recorded for client in row 58 (was 0)
$wsGrupo.Range("I58").Value = 210.15

# Footer counters "X de 58" - one more client now has LAVABOS (I) and
# PORCELANATO (M) sales, so the counts increase by 2 and 1 respectively.
$wsGrupo.Range("I60").Value = "5 de 58"
$wsGrupo.Range("M60").Value = "9 de 58"

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F4").Value = 1405.49
$wsMensual.Range("F15").Value = 7258.4
$wsMensual.Range("F29").Value = 366.83
$wsMensual.Range("F58").Value = 651.2
$wsMensual.Range("F60").Value = 57255.98

# ---------------------------------------------------------------------
# Sheet 3: "CUMPLIMIENTO MENSUAL"
# ---------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 7: LAVABOS
$wsCumpl.Range("D7").Value = 1366.95
$wsCumpl.Range("E7").Value = -480.238983712426
$wsCumpl.Range("F7").Value = 1.541595824221357

# Row 11: PIEDRA SINTERIZADA
$wsCumpl.Range("D11").Value = 13967.89
$wsCumpl.Range("E11").Value = 5605.170249249699
$wsCumpl.Range("F11").Value = 0.7136283147411983

# Row 12: PORCELANATO
$wsCumpl.Range("D12").Value = 27097.3
$wsCumpl.Range("E12").Value = 21526.76
$wsCumpl.Range("F12").Value = 0.5572817243150819

# Row 14: TOTAL
$wsCumpl.Range("D14").Value = 62860.14000000001
$wsCumpl.Range("E14").Value = 37037.85284188786
$wsCumpl.Range("F14").Value = 0.6292432731805835
